# Insert a new price record for "Feria Lagunitas de Puerto Montt" (Mango)
# as row 187, pushing the existing rows 187-241 down to 188-242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 187; Excel shifts rows 187:241 down to 188:242
# and carries the dimension/used-range out to row 242 automatically.
$ws.Rows(187).Insert()

# New record for the inserted row (same market/product template as its
# neighbours, with the new date/quality/price/origin data).
$newRow = @(
    4,                                    # A Mercado ID
    "Feria Lagunitas de Puerto Montt",    # B Mercado
    "Los Lagos",                          # C Región
    44809,                                # D Fecha
    10,                                   # E Codreg
    "Fruta",                              # F Tipo
    100108,                               # G Producto ID
    "Tropicales y subtropicales",         # H Producto
    100108002,                            # I Categoría ID
    "Mango",                              # J Categoría
    "Sin especificar",                    # K Variedad
    "Primera",                            # L Calidad
    60,                                   # M Volumen
    11000,                                # N Precio mínimo
    12000,                                # O Precio máximo
    11500,                                # P Precio promedio ponderado
    '$/bandeja 4 kilos',                  # Q Unidad de comercialización
    "Brasil",                             # R Origen
    2875,                                 # S Precio $/Kg
    4                                     # T Kg / unidad
)

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(187, $i + 1).Value = $newRow[$i]
}
